$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# --- Column B (Coin name) swaps ---
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("B22").Value = 'Dai'
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("B49").Value = 'Maker'

# --- Column C (Link) swaps ---
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'

# --- Column D (Price) updates ---
Set-TextValue "D2" '30.083.65'
Set-TextValue "D3" '1.902.67'
Set-TextValue "D4" '0.9993'
Set-TextValue "D5" '0.8369'
Set-TextValue "D7" '0.9994'
Set-TextValue "D9" '26.70'
Set-TextValue "D10" '0.07072'
Set-TextValue "D11" '0.08088'
Set-TextValue "D12" '0.7651'
Set-TextValue "D13" '1.920.32'
Set-TextValue "D14" '5.275'
Set-TextValue "D15" '92.56'
Set-TextValue "D16" '30.065.83'
Set-TextValue "D17" '14.17'
Set-TextValue "D18" '5.875'
Set-TextValue "D19" '244.76'
Set-TextValue "D20" '0.000007774'
Set-TextValue "D21" '2.158.63'
Set-TextValue "D22" '1.000'
Set-TextValue "D24" '7.024'
Set-TextValue "D25" '0.1763'
Set-TextValue "D26" '9.292'
Set-TextValue "D27" '165.77'
Set-TextValue "D28" '18.97'
Set-TextValue "D29" '2.098'
Set-TextValue "D31" '1.519'
Set-TextValue "D32" '0.05944'
Set-TextValue "D33" '4.298'
Set-TextValue "D34" '4.080'
Set-TextValue "D35" '1.272'
Set-TextValue "D36" '0.7331'
Set-TextValue "D38" '0.01923'
Set-TextValue "D39" '2.782'
Set-TextValue "D40" '0.4454'
Set-TextValue "D41" '73.01'
Set-TextValue "D42" '5.964'
Set-TextValue "D43" '0.8580'
Set-TextValue "D44" '1.909'
Set-TextValue "D45" '0.9988'
Set-TextValue "D46" '102.01'
Set-TextValue "D47" '7.572'
Set-TextValue "D48" '9.843'
Set-TextValue "D49" '1.002.49'
Set-TextValue "D50" '2.060.46'

# --- Column E (Volume 1h) updates ---
$ws.Range("E2").Value = '  -0.01%  '
$ws.Range("E3").Value = '  -0.58%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("E5").Value = '  +4.69%  '
$ws.Range("E6").Value = '  -0.62%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("E8").Value = '  +2.98%  '
$ws.Range("E9").Value = '  +1.17%  '
$ws.Range("E10").Value = '  +1.59%  '
$ws.Range("E11").Value = '  +1.06%  '
$ws.Range("E12").Value = '  +1.69%  '
$ws.Range("E13").Value = '  +0.39%  '
$ws.Range("E14").Value = '  +0.53%  '
$ws.Range("E15").Value = '  -1.17%  '
$ws.Range("E16").Value = '  -0.14%  '
$ws.Range("E17").Value = '  +0.65%  '
$ws.Range("E18").Value = '  -1.81%  '
$ws.Range("E19").Value = '  -1.85%  '
$ws.Range("E20").Value = '  -0.63%  '
$ws.Range("E21").Value = '  -0.19%  '
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("E24").Value = '  +1.08%  '
$ws.Range("E25").Value = '  +24.81%  '
$ws.Range("E26").Value = '  -0.48%  '
$ws.Range("E27").Value = '  -2.11%  '
$ws.Range("E28").Value = '  +0.02%  '
$ws.Range("E29").Value = '  +1.63%  '
$ws.Range("E31").Value = '  -0.71%  '
$ws.Range("E32").Value = '  +8.54%  '
$ws.Range("E33").Value = '  -1.36%  '
$ws.Range("E34").Value = '  -1.19%  '
$ws.Range("E35").Value = '  +0.36%  '
$ws.Range("E36").Value = '  -0.91%  '
$ws.Range("E37").Value = '  -0.44%  '
$ws.Range("E38").Value = '  -0.54%  '
$ws.Range("E39").Value = '  -0.48%  '
$ws.Range("E40").Value = '  -0.18%  '
$ws.Range("E41").Value = '  -0.50%  '
$ws.Range("E42").Value = '  -3.70%  '
$ws.Range("E43").Value = '  +2.75%  '
$ws.Range("E44").Value = '  -0.06%  '
$ws.Range("E45").Value = '  -0.16%  '
$ws.Range("E46").Value = '  +1.22%  '
$ws.Range("E47").Value = '  -0.64%  '
$ws.Range("E48").Value = '  -0.26%  '
$ws.Range("E49").Value = '  +1.39%  '
$ws.Range("E50").Value = '  -0.19%  '
$ws.Range("E51").Value = '  +0.62%  '
